$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C2:C17) from 2023-10-13 (45212)
# to 2023-10-22 (45221) for all data rows.
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
